# Updates cryptocurrency price/volume data in the active worksheet
# to reflect refreshed values from the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.299.04"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "3.640.56"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'197.90"
$ws.Range("E5").Value = "  +11.16%  "
$ws.Range("D6").Value = "'579.56"
$ws.Range("D7").Value = "3.633.06"
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("E11").Value = "  +9.19%  "
$ws.Range("D12").Value = "'56.67"
$ws.Range("E12").Value = "  +7.09%  "
$ws.Range("D13").Value = "'0.0000294"
$ws.Range("E13").Value = "  +18.55%  "
$ws.Range("D14").Value = "'10.11"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").Value = "4.223.01"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "3.642.00"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  +5.08%  "
$ws.Range("D19").Value = "68.302.28"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").Value = "'18.64"
$ws.Range("E20").Value = "  +3.36%  "
$ws.Range("E21").Value = "  +4.68%  "
$ws.Range("D22").Value = "'404.07"
$ws.Range("E22").Value = "  +4.73%  "
$ws.Range("D23").Value = "'13.16"
$ws.Range("E23").Value = "  +30.51%  "
$ws.Range("D24").Value = "'4.27"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "'86.11"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "'2.97"
$ws.Range("E26").Value = "  +5.08%  "
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("D28").Value = "'3.87"
$ws.Range("E28").Value = "  +7.98%  "
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("E30").Value = "  +24.13%  "
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("D33").Value = "'685.05"
$ws.Range("E33").Value = "  +16.75%  "
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  +6.46%  "
$ws.Range("D36").Value = "'64.74"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "'42.79"
$ws.Range("E37").Value = "  +5.68%  "
$ws.Range("E38").Value = "  +17.21%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "0.0₃0792"
$ws.Range("E40").Value = "  +9.36%  "
$ws.Range("D41").Value = "'2.91"
$ws.Range("E41").Value = "  +23.46%  "
$ws.Range("E42").Value = "  +7.66%  "
$ws.Range("E43").Value = "  +16.36%  "
$ws.Range("D44").Value = "3.227.94"
$ws.Range("E44").Value = "  +17.96%  "
$ws.Range("D45").Value = "'3.03"
$ws.Range("E45").Value = "  +41.49%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("D48").Value = "'8.98"
$ws.Range("E48").Value = "  +10.81%  "
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("D50").Value = "'3.14"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("E51").Value = "  +4.10%  "
